$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98 (shifts old rows 98..123 down to 99..124)
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly price record
$ws.Range("A98").Value = 3
$ws.Range("B98").Value = "Femacal de La Calera"
$ws.Range("C98").Value = "Coquimbo"
$ws.Range("D98").Value = "2021-12-29"
$ws.Range("E98").Value = 5
$ws.Range("F98").Value = 100112026
$ws.Range("G98").Value = "Haba"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 65
$ws.Range("K98").Value = 8000
$ws.Range("L98").Value = 8000
$ws.Range("M98").Value = 8000
$ws.Range("N98").Value = "$/saco 25 kilos"
$ws.Range("O98").Value = "Provincia de Petorca"
$ws.Range("P98").Value = 320
$ws.Range("Q98").Value = 25
$ws.Range("R98").Value = "Hortaliza"
